$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.09
$wsSummary.Range("B4").Value = 0.08
$wsSummary.Range("B5").Value = 0.04
$wsSummary.Range("B6").Value = 45
$wsSummary.Range("B8").Value = 18
$wsSummary.Range("B9").Value = 40

# --- Strategy Status sheet (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.09
$wsStatus.Range("D4").Value = 45
$wsStatus.Range("E4").Value = 0.08
$wsStatus.Range("F4").Value = 0.09
$wsStatus.Range("G4").Value = 40

# --- All Trades & MarketMaking sheets: trade #45 (row 46) closed ---
$tradeSheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("G46").Value = 0.91
    $ws.Range("H46").Value = "CLOSED"
    $ws.Range("I46").Value = -2.7752
    $ws.Range("J46").Value = -0.03
    $ws.Range("K46").Value = 100.09
    $ws.Range("P46").Value = "early_exit"
    $ws.Range("Q46").Value = 0.14
}
